$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.241.29"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.215.32"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.94"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.97"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.516"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.63"
$ws.Range("E10").Value = "  +7.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.03"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0785"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.42"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "2.555.79"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.87"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "2.205.10"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.738"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "40.134.30"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.39"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.82"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.83"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.51"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.36"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.93"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.41"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0717"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.98"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.56"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "2.069.03"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.27"
$ws.Range("E44").Value = "  +6.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0271"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.04"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("E47").Value = "  +6.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  -13.33%  "
$ws.Range("D49").Value = "2.427.54"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("E51").Value = "  +1.28%  "
